$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows (comuni MO - rolling 7gg) for 22/04 - 26/04/2021
$data = @(
    @(44308,24,4,6,77,65,16,18,7,12,0,9,53,43,1,1,0,37,8,12,37,255,11,5,16,43,3,11,15,7,40,4,125,11,31,6,26,23,49,5,28,1160,4,0,3,1,4,1,3,0),
    @(44309,23,4,5,62,69,14,14,7,10,1,10,47,42,1,1,0,48,7,8,32,260,7,9,22,43,3,9,10,7,34,5,122,9,25,6,24,22,51,8,35,1134,5,0,3,1,4,2,3,0),
    @(44310,23,6,6,62,67,7,17,7,8,1,9,49,41,1,1,0,53,8,9,25,248,5,8,21,47,3,13,8,8,31,4,112,9,19,6,23,18,50,8,31,1090,5,0,2,1,4,2,4,0),
    @(44311,26,7,11,70,58,6,16,5,8,1,8,41,37,1,1,0,44,6,7,25,252,6,10,19,50,3,10,6,5,30,6,96,11,17,5,22,21,60,7,28,1062,10,0,0,1,3,2,4,0),
    @(44312,30,7,10,67,59,7,18,5,13,1,9,35,39,1,1,0,46,4,7,24,254,9,12,18,49,3,10,5,6,21,7,96,11,17,5,17,22,61,7,29,1062,10,0,1,0,3,2,4,0)
)

$startRow = 234
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $data[$i]

    # Copy style (date column formatting) from the last existing data row
    $ws.Range("A233").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null

    $ws.Range("A$r").Value = $rowValues[0]

    for ($c = 1; $c -lt $rowValues.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c + 1)
        $cell.Value = $rowValues[$c]
    }
}

$excel.CutCopyMode = 0
